$wb = $excel.ActiveWorkbook
$wsGit = $wb.Worksheets.Item("gitHUB")
$wsTest = $wb.Worksheets.Item("myTest")

# --- myTest sheet: insert two new columns (F,G) for Authentication fields ---
$wsTest.Range("F1:G1").EntireColumn.Insert()
$wsTest.Range("F1:G1").ColumnWidth = 18.17

# New header row values
$wsTest.Range("F1").Value = "AuthenticationType"
$wsTest.Range("G1").Value = "AuthenticationDetails"

# Replace old demo title/brand values (previously at F/G, now shifted to H/I)
# with the new "ourTest" title value; brand keeps the plain "brand" text.
$wsTest.Range("H2").Value = "ourTest"
$wsTest.Range("I2").Value = "brand"
$wsTest.Range("H4").Value = "ourTest"
$wsTest.Range("I4").Value = "brand"

# Update selection/active sheet state: myTest becomes the active/selected sheet
$wsGit.Range("F1").Select()
$wsTest.Range("C4").Select()
$wsTest.Activate()

Write-Host "edit complete"
